$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Hora") goes from 6 to 7 for every data row (2-51)
$ws.Range("G2:G51").Value = "'7"

# Column D ("Price") updates for the rows where the scraped price changed
$ws.Range("D2").Value = "'260.78"
$ws.Range("D3").Value = "'21.55"
$ws.Range("D4").Value = "'6.237"
$ws.Range("D5").Value = "'0.06157"
$ws.Range("D6").Value = "'3.576"
$ws.Range("D7").Value = "'6.552"
$ws.Range("D8").Value = "'1.406"
$ws.Range("D9").Value = "'0.8235"
$ws.Range("D11").Value = "'0.08185"
$ws.Range("D12").Value = "'0.03549"
$ws.Range("D13").Value = "'0.03189"
$ws.Range("D14").Value = "'0.09202"
$ws.Range("D15").Value = "'3.773"
$ws.Range("D16").Value = "'0.001630"
$ws.Range("D17").Value = "'0.04660"
$ws.Range("D18").Value = "'0.006442"
$ws.Range("D19").Value = "'0.006164"
$ws.Range("D20").Value = "'0.001071"
$ws.Range("D21").Value = "'0.0001502"
$ws.Range("D22").Value = "'3.726"
$ws.Range("D23").Value = "'2.269"
$ws.Range("D25").Value = "'0.3315"
$ws.Range("D40").Value = "'0.04677"
$ws.Range("D41").Value = "'0.006989"
$ws.Range("D42").Value = "'0.003763"
$ws.Range("D44").Value = "'0.01178"
$ws.Range("D45").Value = "'0.00006084"
$ws.Range("D48").Value = "'0.9820"
$ws.Range("D49").Value = "'0.001135"
$ws.Range("D50").Value = "'0.00001903"
$ws.Range("D51").Value = "'0.01242"
